# Refactor code structure + fix pip package issue (#9)
#
# Updates the test-result workbook to reflect a run where the routing-table
# checks failed (wrong route counts) instead of passing, and removes the
# now-unreached "pc.net reachable" checks from the "All" sheet. The
# "Failed" sheet is populated with the 9 newly-failing checks, and the
# "Summary" sheet counters are refreshed.

$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("Summary")
$wsAll     = $wb.Worksheets.Item("All")
$wsFailed  = $wb.Worksheets.Item("Failed")

# The 9 routing-table checks that now fail, in sheet order. Each entry is
# (description-suffix-used-on-"Failed"-sheet, reason).
$failures = @(
    @("as1r1", "The routing table of as1r1 have the wrong number of routes: 6, expected: 8"),
    @("as1r2", "The routing table of as1r2 have the wrong number of routes: 5, expected: 8"),
    @("as2r1", "The routing table of as2r1 have the wrong number of routes: 5, expected: 8"),
    @("as2r2", "The routing table of as2r2 have the wrong number of routes: 4, expected: 8"),
    @("as3r1", "The routing table of as3r1 have the wrong number of routes: 7, expected: 9"),
    @("root",  "The routing table of root have the wrong number of routes: 1, expected: 2"),
    @("net",   "The routing table of net have the wrong number of routes: 1, expected: 2"),
    @("pc",    "The routing table of pc have the wrong number of routes: 1, expected: 2"),
    @("local", "The routing table of local have the wrong number of routes: 1, expected: 2")
)

# --- "All" sheet: rows 332..340 flip from True/OK to False/<reason> -------
$firstRow = 332
for ($i = 0; $i -lt $failures.Length; $i++) {
    $row = $firstRow + $i
    $reason = $failures[$i][1]
    $wsAll.Range("B$row").Value = "'False"
    $wsAll.Range("C$row").Value = $reason
}

# --- "All" sheet: drop the 6 now-removed `pc.net` reachability rows -------
$wsAll.Range("A352:A357").EntireRow.Delete()

# --- "Failed" sheet: append the 9 failing checks (rows 2..10) -------------
$failedRow = 2
for ($i = 0; $i -lt $failures.Length; $i++) {
    $device = $failures[$i][0]
    $reason = $failures[$i][1]
    $wsFailed.Range("A$failedRow").Value = "Checking the routing table of $device"
    $wsFailed.Range("B$failedRow").Value = "'False"
    $wsFailed.Range("C$failedRow").Value = $reason
    $failedRow = $failedRow + 1
}

# --- "Summary" sheet: refresh the totals ----------------------------------
$wsSummary.Range("A2").Value = "'350"
$wsSummary.Range("B2").Value = "'341"
$wsSummary.Range("C2").Value = "'9"
